$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"

$ws.Range("B7").Value = 0.04366584981362852
$ws.Range("C7").Value = 0.2151581326876664
$ws.Range("D7").Value = 0.06088480029806752
$ws.Range("E7").Value = 0.2467484555130336
$ws.Range("F7").Value = 0.2575856276859916
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.05516078357601667
$ws.Range("C8").Value = 0.3645461267824502
$ws.Range("D8").Value = 0.1795074568246066
$ws.Range("E8").Value = 0.4236832033779562
$ws.Range("F8").Value = 0.4601713743116403
$ws.Range("G8").Value = 6

$ws.Range("B9").Value = 0.1972001695263083
$ws.Range("C9").Value = 0.2323019497277908
$ws.Range("D9").Value = 0.138851959107692
$ws.Range("E9").Value = 0.3726284464553022
$ws.Range("F9").Value = 0.3872287158382381
$ws.Range("G9").Value = 3

$ws.Range("B10").Value = -0.06854498788710228
$ws.Range("C10").Value = 0.06854498788710228
$ws.Range("D10").Value = 0.004698415364442998
$ws.Range("E10").Value = 0.06854498788710228
$ws.Range("G10").Value = 1
